$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the field name: "qualifiedProfessionalOrganization" -> "qualifiedProfessionalOrganizati"
#    (AGOL's internal field name is truncated to 31 characters)
$ws.Range("A25").Value = "qualifiedProfessionalOrganizati"

# 2. Add reviewer comment on A25 explaining why the field name should not be "fixed"
$comment = $ws.Range("A25").AddComment("Roy Jeong:" + [char]10 + "Never rename fields. The field name is not incorrect, it is currently actually referenced as such in AGOL internally.")

# 3. Update the active selection to B29 (matches the cursor position after the edit)
$ws.Range("B29").Select() | Out-Null
